$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.189921424373097
$ws.Range("D2").Value = 0.1170785113430952
$ws.Range("E2").Value = 0.1342225034446365
$ws.Range("F2").Value = 2.299143291757318
$ws.Range("G2").Value = 0.002544690869663589
$ws.Range("K2").Value = 0.9768027924936291
$ws.Range("L2").Value = 0.2011692857141156
$ws.Range("M2").Value = 0.2709987220115515
$ws.Range("N2").Value = 3.409348728329093
$ws.Range("B3").Value = 1.170942302046143
$ws.Range("D3").Value = 0.1175843156271625
$ws.Range("E3").Value = 0.1320792103635089
$ws.Range("F3").Value = 2.265056052188172
$ws.Range("G3").Value = 0.002549534415284718
$ws.Range("K3").Value = 0.8949374470610394
$ws.Range("L3").Value = 0.1942138136348603
$ws.Range("M3").Value = 0.2651847245926753
$ws.Range("N3").Value = 3.415377617336276
$ws.Range("B4").Value = 1.159897357817812
$ws.Range("D4").Value = 0.1179106796052896
$ws.Range("E4").Value = 0.1307364254393182
$ws.Range("F4").Value = 2.245257393689243
$ws.Range("G4").Value = 0.002552665839556659
$ws.Range("K4").Value = 0.8452283808001653
$ws.Range("L4").Value = 0.19006707922874
$ws.Range("M4").Value = 0.261764430451791
$ws.Range("N4").Value = 3.419835966827264
$ws.Range("B5").Value = 1.155549602475219
$ws.Range("D5").Value = 0.1180476455209352
$ws.Range("E5").Value = 0.1301824095485244
$ws.Range("F5").Value = 2.237472949194014
$ws.Range("G5").Value = 0.002553981648961803
$ws.Range("K5").Value = 0.8251103648763092
$ws.Range("L5").Value = 0.1884083653463478
$ws.Range("M5").Value = 0.2604082434075927
$ws.Range("N5").Value = 3.421842780764919
$ws.Range("B6").Value = 1.154836916367032
$ws.Range("D6").Value = 0.1180706283090274
$ws.Range("E6").Value = 0.1300900013125901
$ws.Range("F6").Value = 2.236197459650398
$ws.Range("G6").Value = 0.002554202541368117
$ws.Range("K6").Value = 0.8217781362181995
$ws.Range("L6").Value = 0.1881348148125284
$ws.Range("M6").Value = 0.2601853215495886
$ws.Range("N6").Value = 3.422187479010518
$ws.Range("B7").Value = 1.159838102129072
$ws.Range("D7").Value = 0.1179125107051355
$ws.Range("E7").Value = 0.1307289815083603
$ws.Range("F7").Value = 2.245151262372445
$ws.Range("G7").Value = 0.002552683424373395
$ws.Range("K7").Value = 0.8449565015527014
$ws.Range("L7").Value = 0.190044583325232
$ws.Range("M7").Value = 0.2617459881391468
$ws.Range("N7").Value = 3.419862262471483
$ws.Range("B8").Value = 1.183251267409901
$ws.Range("D8").Value = 0.1172496321936762
$ws.Range("E8").Value = 0.1334890069282135
$ws.Range("F8").Value = 2.287154799186666
$ws.Range("G8").Value = 0.002546328321201453
$ws.Range("K8").Value = 0.9484593938609009
$ws.Range("L8").Value = 0.1987452761203485
$ws.Range("M8").Value = 0.2689630260727824
$ws.Range("N8").Value = 3.411270306209516
$ws.Range("B9").Value = 1.233987098606548
$ws.Range("D9").Value = 0.1160751844918586
$ws.Range("E9").Value = 0.1386927040949839
$ws.Range("F9").Value = 2.378538442708759
$ws.Range("G9").Value = 0.002535109410193242
$ws.Range("K9").Value = 1.155906865955899
$ws.Range("L9").Value = 0.2167945242097744
$ws.Range("M9").Value = 0.2843026738090728
$ws.Range("N9").Value = 3.400437120397086
$ws.Range("B10").Value = 1.274203383047904
$ws.Range("D10").Value = 0.115288923342149
$ws.Range("E10").Value = 0.1423942765785355
$ws.Range("F10").Value = 2.451240019730619
$ws.Range("G10").Value = 0.002527616482589395
$ws.Range("K10").Value = 1.311151795353283
$ws.Range("L10").Value = 0.2306640182699766
$ws.Range("M10").Value = 0.2962989841750527
$ws.Range("N10").Value = 3.396164173540001
$ws.Range("B11").Value = 1.293137956132881
$ws.Range("D11").Value = 0.1149478975028835
$ws.Range("E11").Value = 0.14405305822849
$ws.Range("F11").Value = 2.485536696830792
$ws.Range("G11").Value = 0.002524368729715178
$ws.Range("K11").Value = 1.382415361999961
$ws.Range("L11").Value = 0.2371073997377806
$ws.Range("M11").Value = 0.3019148095152389
$ws.Range("N11").Value = 3.39502500060081
$ws.Range("B12").Value = 1.300399949649687
$ws.Range("D12").Value = 0.1148211578882083
$ws.Range("E12").Value = 0.1446776851010689
$ws.Range("F12").Value = 2.498701060620562
$ws.Range("G12").Value = 0.002523161878719021
$ws.Range("K12").Value = 1.40949488719906
$ws.Range("L12").Value = 0.2395667183419476
$ws.Range("M12").Value = 0.3040642131454092
$ws.Range("N12").Value = 3.394709663160782
$ws.Range("B13").Value = 1.298831865243045
$ws.Range("D13").Value = 0.1148483467246706
$ws.Range("E13").Value = 0.1445433155523617
$ws.Range("F13").Value = 2.495857994753351
$ws.Range("G13").Value = 0.002523420774710118
$ws.Range("K13").Value = 1.403658645614371
$ws.Range("L13").Value = 0.2390361983609353
$ws.Range("M13").Value = 0.3036002860148415
$ws.Range("N13").Value = 3.394772410206386
$ws.Range("B14").Value = 1.293733564091156
$ws.Range("D14").Value = 0.1149374224376523
$ws.Range("E14").Value = 0.1441045166441874
$ws.Range("F14").Value = 2.486616185318951
$ws.Range("G14").Value = 0.002524268981091439
$ws.Range("K14").Value = 1.384641328491341
$ws.Range("L14").Value = 0.2373093411549547
$ws.Range("M14").Value = 0.3020911850464856
$ws.Range("N14").Value = 3.39499672996439
$ws.Range("B15").Value = 1.290622668839717
$ws.Range("D15").Value = 0.1149922965323356
$ws.Range("E15").Value = 0.1438352841898265
$ws.Range("F15").Value = 2.48097838192939
$ws.Range("G15").Value = 0.002524791524136795
$ws.Range("K15").Value = 1.373004894492453
$ws.Range("L15").Value = 0.2362541132729916
$ws.Range("M15").Value = 0.3011697879999318
$ws.Range("N15").Value = 3.395149254663934
$ws.Range("B16").Value = 1.272978829574157
$ws.Range("D16").Value = 0.1153115456495115
$ws.Range("E16").Value = 0.1422853741686083
$ws.Range("F16").Value = 2.449023364705681
$ws.Range("G16").Value = 0.002527831956374576
$ws.Range("K16").Value = 1.306507576933257
$ws.Range("L16").Value = 0.2302456320492325
$ws.Range("M16").Value = 0.2959351690371221
$ws.Range("N16").Value = 3.396254840990537
$ws.Range("B17").Value = 1.262318714485986
$ws.Range("D17").Value = 0.1155116624743151
$ws.Range("E17").Value = 0.1413281934693646
$ws.Range("F17").Value = 2.429734217792458
$ws.Range("G17").Value = 0.002529738264852992
$ws.Range("K17").Value = 1.265878914584789
$ws.Range("L17").Value = 0.2265940131647852
$ws.Range("M17").Value = 0.2927645316007315
$ws.Range("N17").Value = 3.397139403233055
$ws.Range("B18").Value = 1.256247545027662
$ws.Range("D18").Value = 0.1156283310526973
$ws.Range("E18").Value = 0.1407752822032275
$ws.Range("F18").Value = 2.418754799404326
$ws.Range("G18").Value = 0.002530849867085608
$ws.Range("K18").Value = 1.242570721971617
$ws.Range("L18").Value = 0.2245063144577415
$ws.Range("M18").Value = 0.2909557970226615
$ws.Range("N18").Value = 3.39772389325951
$ws.Range("B19").Value = 1.254202303373063
$ws.Range("D19").Value = 0.1156681019436796
$ws.Range("E19").Value = 0.1405876668476296
$ws.Range("F19").Value = 2.415057113732274
$ws.Range("G19").Value = 0.002531228841180971
$ws.Range("K19").Value = 1.234689293733879
$ws.Range("L19").Value = 0.2238016202507538
$ws.Range("M19").Value = 0.2903459552273375
$ws.Range("N19").Value = 3.397934785061622
$ws.Range("B20").Value = 1.263447267712792
$ws.Range("D20").Value = 0.1154901975138714
$ws.Range("E20").Value = 0.1414303312402936
$ws.Range("F20").Value = 2.43177565241777
$ws.Range("G20").Value = 0.002529533768597935
$ws.Range("K20").Value = 1.270197651085596
$ws.Range("L20").Value = 0.2269814280344207
$ws.Range("M20").Value = 0.2931005059496599
$ws.Range("N20").Value = 3.397037401641214
$ws.Range("B21").Value = 1.295228566231316
$ws.Range("D21").Value = 0.1149111935765923
$ws.Range("E21").Value = 0.14423349727603
$ws.Range("F21").Value = 2.489325919027493
$ws.Range("G21").Value = 0.002524019218948239
$ws.Range("K21").Value = 1.3902246262208
$ws.Range("L21").Value = 0.237816035263009
$ws.Range("M21").Value = 0.3025338254536507
$ws.Range("N21").Value = 3.394927689707984
$ws.Range("B22").Value = 1.316534973661817
$ws.Range("D22").Value = 0.1145467696973483
$ws.Range("E22").Value = 0.1460450662916557
$ws.Range("F22").Value = 2.527970208222598
$ws.Range("G22").Value = 0.002520549161706307
$ws.Range("K22").Value = 1.469215125112328
$ws.Range("L22").Value = 0.2450098968630101
$ws.Range("M22").Value = 0.3088320144231105
$ws.Range("N22").Value = 3.394225418250301
$ws.Range("B23").Value = 1.305114392858513
$ws.Range("D23").Value = 0.1147399877405419
$ws.Range("E23").Value = 0.1450800404104449
$ws.Range("F23").Value = 2.50725030905997
$ws.Range("G23").Value = 0.002522388974143132
$ws.Range("K23").Value = 1.427006035469446
$ws.Range("L23").Value = 0.2411600511252772
$ws.Range("M23").Value = 0.3054583854466131
$ws.Range("N23").Value = 3.394538215667268
$ws.Range("B24").Value = 1.262936869768509
$ws.Range("D24").Value = 0.1154998967802001
$ws.Range("E24").Value = 0.1413841629090697
$ws.Range("F24").Value = 2.43085237671977
$ws.Range("G24").Value = 0.002529626172482576
$ws.Range("K24").Value = 1.26824499540021
$ws.Range("L24").Value = 0.2268062414411247
$ws.Range("M24").Value = 0.2929485680050377
$ws.Range("N24").Value = 3.397083280015309
$ws.Range("B25").Value = 1.219745363252144
$ws.Range("D25").Value = 0.1163794609581466
$ws.Range("E25").Value = 0.1373067288337122
$ws.Range("F25").Value = 2.352844838486533
$ws.Range("G25").Value = 0.002538012178442696
$ws.Range("K25").Value = 1.099296781303536
$ws.Range("L25").Value = 0.2118053314559063
$ws.Range("M25").Value = 0.2800255455665663
$ws.Range("N25").Value = 3.402721836176283
